# Append four new DBC message blocks (DYN_F_SIG1, DYN_F_SIG2, DYN_R_SIG1,
# DYN_R_SIG2) to the "Autonomous_temporary" sheet, starting at row 42,
# mirroring the formatting of the existing message blocks on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autonomous_temporary")

$headerFillBlue = 15652797   # BGR for FGColor/BgColor 00BDD7EE
$headerFillGold = 6740479    # BGR for FGColor/BgColor 00FFD966

$colHeaders = @("Signal Name", "Start Bit", "Length (bits)", "Byte Order", "Signed", "Factor", "Offset", "Min", "Max", "Unit", "Choices")

function Write-MessageHeader($row, $msgName, $msgId, $sender) {
    $ws.Cells.Item($row, 1).Value = "Message: $msgName"
    $ws.Cells.Item($row, 2).Value = "ID: $msgId"
    $ws.Cells.Item($row, 3).Value = "Sender(s): $sender"
    $hdrRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 3))
    $hdrRange.Font.Bold = $true
    $hdrRange.Interior.Color = $headerFillBlue
    $hdrRange.Borders.LineStyle = 1
}

function Write-ColumnHeader($row) {
    for ($c = 1; $c -le 11; $c++) {
        $ws.Cells.Item($row, $c).Value = $colHeaders[$c - 1]
    }
    $hdrRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 11))
    $hdrRange.Font.Bold = $true
    $hdrRange.Interior.Color = $headerFillGold
    $hdrRange.Borders.LineStyle = 1
}

function Write-SignalRow($row, $name, $startBit, $length, $byteOrder, $signed, $factor, $offset, $min, $max, $unit, $choices) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $startBit
    $ws.Cells.Item($row, 3).Value = $length
    $ws.Cells.Item($row, 4).Value = $byteOrder
    $ws.Cells.Item($row, 5).Value = $signed
    $ws.Cells.Item($row, 6).Value = $factor
    $ws.Cells.Item($row, 7).Value = $offset
    if ($null -ne $min) { $ws.Cells.Item($row, 8).Value = $min }
    if ($null -ne $max) { $ws.Cells.Item($row, 9).Value = $max }
    if ($unit) { $ws.Cells.Item($row, 10).Value = $unit }
    if ($choices) { $ws.Cells.Item($row, 11).Value = $choices }
    $dataRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 11))
    $dataRange.Borders.LineStyle = 1
}

function Touch-BlankRow($row) {
    # Forces the worksheet to materialise an explicit (but empty) <row>
    # element, matching the separator rows already used between message
    # blocks elsewhere in this sheet -- without changing any cell content.
    $r = $ws.Rows.Item($row)
    $r.OutlineLevel = 0
}

# Row 42 left blank (separator row).
Touch-BlankRow 42

# --- Message: DYN_F_SIG1 ---
Write-MessageHeader 43 "DYN_F_SIG1" "0x1be" "DYN_F"
Write-ColumnHeader 44
Write-SignalRow 45 "ST_ANGLE" 0  16 "Intel" $true  10 0 $null $null "º"    $null
Write-SignalRow 46 "SUSP_R"   16 16 "Intel" $false 10 0 $null $null "mm"   $null
Write-SignalRow 47 "SUSP_L"   32 16 "Intel" $false 10 0 $null $null "mm"   $null

# Row 48 left blank (separator row).
Touch-BlankRow 48

# --- Message: DYN_F_SIG2 ---
Write-MessageHeader 49 "DYN_F_SIG2" "0x1c8" "DYN_F"
Write-ColumnHeader 50
Write-SignalRow 51 "SPD_LEFT"  0  16 "Intel" $false 10 0 $null $null "km/h" $null
Write-SignalRow 52 "SPD_RIGHT" 16 16 "Intel" $false 10 0 $null $null "km/h" $null

# Row 53 left blank (separator row).
Touch-BlankRow 53

# --- Message: DYN_R_SIG1 ---
Write-MessageHeader 54 "DYN_R_SIG1" "0x222" "DYN_R"
Write-ColumnHeader 55
Write-SignalRow 56 "BRK_PRESS" 0  16 "Intel" $false 10 0 $null $null "bar" $null
Write-SignalRow 57 "SUSP_R"    16 16 "Intel" $false 10 0 $null $null "mm"  $null
Write-SignalRow 58 "SUSP_L"    32 16 "Intel" $false 10 0 $null $null "mm"  $null

# Row 59 left blank (separator row).
Touch-BlankRow 59

# --- Message: DYN_R_SIG2 ---
Write-MessageHeader 60 "DYN_R_SIG2" "0x22c" "DYN_R"
Write-ColumnHeader 61
Write-SignalRow 62 "SPD_LEFT"  0  16 "Intel" $false 10 0 $null $null "km/h" $null
Write-SignalRow 63 "SPD_RIGHT" 16 16 "Intel" $false 10 0 $null $null "km/h" $null
